$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new text could be misread as a number by Excel's
# automatic type inference are pre-formatted as Text so the exact
# original string (including trailing zeros) survives the round-trip.

$ws.Range("D2").Value = "59.417.56"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "2.630.13"
$ws.Range("E3").Value = "  +3.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.57"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.74"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "2.655.09"
$ws.Range("E9").Value = "  +4.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.29"
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("E11").Value = "  +3.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("D14").Value = "3.093.76"
$ws.Range("E14").Value = "  +3.63%  "
$ws.Range("D15").Value = "59.279.56"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.07"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "2.649.11"
$ws.Range("E18").Value = "  +4.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "349.90"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.53"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("E22").Value = "  +4.60%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.78"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.420"
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.164"
$ws.Range("E26").Value = "  +3.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").Value = "0.0₃0812"
$ws.Range("E28").Value = "  +4.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.14"
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("E31").Value = "  +8.63%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.59"
$ws.Range("E32").Value = "  +3.91%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.01"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.10"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.984"
$ws.Range("E35").Value = "  +9.69%  "
$ws.Range("E36").Value = "  +4.07%  "
$ws.Range("E37").Value = "  +3.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.79"
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.849"
$ws.Range("E39").Value = "  +4.31%  "
$ws.Range("E40").Value = "  +5.54%  "
$ws.Range("E41").Value = "  +3.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "278.71"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0988"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.610"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.71"
$ws.Range("E46").Value = "  +6.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0525"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0231"
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.74"
$ws.Range("E49").Value = "  +5.29%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.29"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.987.63"
$ws.Range("E51").Value = "  +5.23%  "
